$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.057.62"
$ws.Range("E2").Value = "  +2.05%  "

$ws.Range("D3").Value = "3.095.31"
$ws.Range("E3").Value = "  +4.86%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.13%  "

$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").Value = "3.094.28"
$ws.Range("E8").Value = "  +4.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.483"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.66%  "

$ws.Range("E15").Value = "  -0.18%  "

$ws.Range("D16").Value = "3.602.08"
$ws.Range("E16").Value = "  +4.70%  "

$ws.Range("D17").Value = "66.951.69"
$ws.Range("E17").Value = "  +1.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.57%  "

$ws.Range("D19").Value = "3.093.28"
$ws.Range("E19").Value = "  +4.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +17.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "470.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.714"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.81%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.73%  "

$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.60%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.97%  "

$ws.Range("E32").Value = "  +4.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.14%  "

$ws.Range("E34").Value = "  +6.25%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("E36").Value = "  +3.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.314"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.76%  "

$ws.Range("E43").Value = "  +3.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.48%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "392.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.59%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0365"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.82%  "

$ws.Range("D47").Value = "2.763.83"
$ws.Range("E47").Value = "  +2.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.72%  "

$ws.Range("E49").Value = "  +0.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.86%  "
